$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old placeholder "null" entries that used to mark empty
# reads from the sheet.
$ws.Range("A5").Clear()
$ws.Range("A6").Clear()
$ws.Range("B6").Clear()
$ws.Range("B7").Clear()

# A7 used to hold the placeholder "null" value - it now holds the real
# value that was actually read back ("Admin").
$ws.Range("A7").Value = "Admin"

# B8 used to duplicate A8 ("admin"); now it stores the literal value that
# was read back for an empty cell.
$ws.Range("B8").Value = "dkjhdbsdsk"

# Update the view so the selection reflects the edited cell (this also
# resets the sheet's top-left scroll position back to A1).
$ws.Range("A7").Select()
